# Updates the cryptos price/volume table with freshly scraped values.
# Values are written as text (matching the original inlineStr cells), so a
# temporary "@" (text) number format is applied before the write and then
# cleared again afterwards to avoid leaving any stray cell formatting behind
# (this also prevents Excel from "helpfully" re-interpreting strings such as
# "0.9980" or "29.879.98" as numbers/dates).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ref, $val) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextCell "D2" "29.879.98"
Set-TextCell "E2" "  -1.32%  "
Set-TextCell "D3" "1.894.54"
Set-TextCell "E3" "  -1.37%  "
Set-TextCell "E4" "  -0.04%  "
Set-TextCell "D5" "0.7758"
Set-TextCell "E5" "  -4.18%  "
Set-TextCell "D6" "244.80"
Set-TextCell "E6" "  +0.12%  "
Set-TextCell "D7" "1.000"
Set-TextCell "E7" "  -0.03%  "
Set-TextCell "D8" "0.3148"
Set-TextCell "E8" "  -3.23%  "
Set-TextCell "D9" "0.07557"
Set-TextCell "E9" "  +4.13%  "
Set-TextCell "D10" "25.53"
Set-TextCell "E10" "  -5.43%  "
Set-TextCell "D11" "0.08113"
Set-TextCell "E11" "  +0.28%  "
Set-TextCell "D12" "0.7711"
Set-TextCell "E12" "  -2.47%  "
Set-TextCell "D13" "5.482"
Set-TextCell "E13" "  +1.32%  "
Set-TextCell "D14" "1.863.32"
Set-TextCell "E14" "  -2.97%  "
Set-TextCell "D15" "92.36"
Set-TextCell "E15" "  -1.73%  "
Set-TextCell "D16" "6.228"
Set-TextCell "E16" "  +2.58%  "
Set-TextCell "D17" "29.783.38"
Set-TextCell "D18" "14.03"
Set-TextCell "E18" "  -1.52%  "
Set-TextCell "B19" "BitcoinCash"
Set-TextCell "C19" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextCell "D19" "244.65"
Set-TextCell "E19" "  -2.21%  "
Set-TextCell "B20" "ShibaInu"
Set-TextCell "C20" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextCell "D20" "0.000007916"
Set-TextCell "E20" "  +0.87%  "
Set-TextCell "D21" "0.9995"
Set-TextCell "E22" "  -1.18%  "
Set-TextCell "D23" "2.109.78"
Set-TextCell "E23" "  -3.47%  "
Set-TextCell "E24" "  -0.08%  "
Set-TextCell "E25" "  -5.65%  "
Set-TextCell "D26" "9.463"
Set-TextCell "E26" "  -0.37%  "
Set-TextCell "D27" "162.92"
Set-TextCell "E27" "  -3.07%  "
Set-TextCell "D28" "18.82"
Set-TextCell "E28" "  -1.09%  "
Set-TextCell "E29" "  -5.18%  "
Set-TextCell "D30" "1.435"
Set-TextCell "E30" "  +3.18%  "
Set-TextCell "D31" "1.551"
Set-TextCell "E31" "  +0.05%  "
Set-TextCell "D32" "4.490"
Set-TextCell "E32" "  +3.29%  "
Set-TextCell "E33" "  -1.25%  "
Set-TextCell "D34" "0.05511"
Set-TextCell "E34" "  -4.17%  "
Set-TextCell "D35" "1.262"
Set-TextCell "E35" "  -2.71%  "
Set-TextCell "D36" "0.7583"
Set-TextCell "E36" "  +1.13%  "
Set-TextCell "D37" "0.9980"
Set-TextCell "E37" "  -0.51%  "
Set-TextCell "D38" "2.640"
Set-TextCell "E38" "  -3.21%  "
Set-TextCell "D39" "0.01928"
Set-TextCell "E39" "  -1.61%  "
Set-TextCell "D40" "2.791"
Set-TextCell "E40" "  -1.08%  "
Set-TextCell "D41" "1.162.48"
Set-TextCell "E41" "  +12.45%  "
Set-TextCell "D42" "74.15"
Set-TextCell "E42" "  -0.26%  "
Set-TextCell "D43" "0.4454"
Set-TextCell "E43" "  -2.13%  "
Set-TextCell "D44" "5.944"
Set-TextCell "E44" "  -0.64%  "
Set-TextCell "D45" "0.8477"
Set-TextCell "E45" "  -0.71%  "
Set-TextCell "E46" "  +0.00%  "
Set-TextCell "D47" "1.902"
Set-TextCell "E47" "  -1.42%  "
Set-TextCell "D48" "3.129"
Set-TextCell "E48" "  +0.19%  "
Set-TextCell "D49" "102.23"
Set-TextCell "E49" "  -1.46%  "
Set-TextCell "D50" "9.966"
Set-TextCell "E50" "  -0.29%  "
Set-TextCell "D51" "7.529"
Set-TextCell "E51" "  -1.32%  "
